$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.191.91'
$ws.Range("E2").Value = '  +0.29%  '

$ws.Range("D3").Value = '3.814.66'
$ws.Range("E3").Value = '  -1.19%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").Value = '706.43'
$ws.Range("E5").Value = '  +1.13%  '

$ws.Range("D6").Value = '171.95'
$ws.Range("E6").Value = '  -0.80%  '

$ws.Range("D7").Value = '3.813.87'
$ws.Range("E7").Value = '  -1.17%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").Value = '0.525'
$ws.Range("E9").Value = '  -0.12%  '

$ws.Range("D10").Value = '0.162'
$ws.Range("E10").Value = '  -0.44%  '

$ws.Range("D11").Value = '7.65'
$ws.Range("E11").Value = '  +6.16%  '

$ws.Range("D12").Value = '0.462'
$ws.Range("E12").Value = '  +0.36%  '

$ws.Range("E13").Value = '  -1.64%  '

$ws.Range("D14").Value = '35.98'
$ws.Range("E14").Value = '  -0.84%  '

$ws.Range("D15").Value = '4.458.92'
$ws.Range("E15").Value = '  -1.15%  '

$ws.Range("D16").Value = '3.803.21'
$ws.Range("E16").Value = '  -1.46%  '

$ws.Range("D17").Value = '71.159.73'
$ws.Range("E17").Value = '  +0.14%  '

$ws.Range("D18").Value = '7.17'
$ws.Range("E18").Value = '  -0.51%  '

$ws.Range("D19").Value = '17.42'
$ws.Range("E19").Value = '  -0.06%  '

$ws.Range("E20").Value = '  -0.08%  '

$ws.Range("D21").Value = '511.75'
$ws.Range("E21").Value = '  +2.49%  '

$ws.Range("D22").Value = '10.70'
$ws.Range("E22").Value = '  -0.33%  '

$ws.Range("D23").Value = '0.725'
$ws.Range("E23").Value = '  +0.52%  '

$ws.Range("D24").Value = '84.51'
$ws.Range("E24").Value = '  -0.58%  '

$ws.Range("E25").Value = '  -2.72%  '

$ws.Range("D26").Value = '3.965.13'
$ws.Range("E26").Value = '  -1.24%  '

$ws.Range("D27").Value = '12.06'
$ws.Range("E27").Value = '  -1.35%  '

$ws.Range("D28").Value = '10.41'
$ws.Range("E28").Value = '  -2.66%  '

$ws.Range("E29").Value = '  +0.10%  '

$ws.Range("D30").Value = '2.05'
$ws.Range("E30").Value = '  -3.74%  '

$ws.Range("D31").Value = '3.05'
$ws.Range("E31").Value = '  -3.48%  '

$ws.Range("E32").Value = '  -1.17%  '

$ws.Range("D33").Value = '7.39'
$ws.Range("E33").Value = '  -2.25%  '

$ws.Range("D34").Value = '29.12'
$ws.Range("E34").Value = '  -1.59%  '

$ws.Range("E35").Value = '  -4.81%  '

$ws.Range("D36").Value = '9.17'
$ws.Range("E36").Value = '  -0.41%  '

$ws.Range("D37").Value = '3.778.72'
$ws.Range("E37").Value = '  -1.01%  '

$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.04%  '

$ws.Range("D39").Value = '0.102'
$ws.Range("E39").Value = '  -2.19%  '

$ws.Range("E40").Value = '  -1.13%  '

$ws.Range("E41").Value = '  -2.94%  '

$ws.Range("E42").Value = '  -1.63%  '

$ws.Range("E43").Value = '  -4.49%  '

$ws.Range("E44").Value = '  +0.00%  '

$ws.Range("D46").Value = '167.21'
$ws.Range("E46").Value = '  +1.98%  '

$ws.Range("D47").Value = '0.000313'
$ws.Range("E47").Value = '  -0.49%  '

$ws.Range("D48").Value = '49.18'
$ws.Range("E48").Value = '  +0.00%  '

$ws.Range("D49").Value = '423.84'
$ws.Range("E49").Value = '  +2.12%  '

$ws.Range("D50").Value = '8.63'
$ws.Range("E50").Value = '  -0.21%  '

$ws.Range("E51").Value = '  +4.63%  '
